$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 10 (B10/C10): "Fornecer aos alunos..." -> "7455355 - Robson da Silva Rocha"
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 2).Value = "7455355 - Robson da Silva Rocha"
$ws.Cells.Item(10, 3).Value = "7455355 - Robson da Silva Rocha"

# ---------------------------------------------------------------------------
# 2. Stash copies of the two row "shapes" we need as formatting templates
#    into scratch rows far below the data (rows 100 = ht60 full A/B/C style,
#    101 = ht120 full A/B/C style, 102 = A-only style) before we start
#    clearing/overwriting rows 13-21 (two of which - 16 & 17 - are
#    themselves one of the templates).
# ---------------------------------------------------------------------------
$ws.Range("A10:C10").Copy()
$ws.Range("A100:C100").PasteSpecial(-4122)

$ws.Range("A16:C16").Copy()
$ws.Range("A101:C101").PasteSpecial(-4122)

$ws.Range("A11:C11").Copy()
$ws.Range("A102:C102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Remove the old row 22 (long "Bibliografia" reference text), which is no
#    longer present after the edit (sheet shrinks from A1:C22 to A1:C21).
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).Delete()

# ---------------------------------------------------------------------------
# 4. Rewrite rows 13-21 with their new contents / formatting / heights using
#    the stashed templates.
# ---------------------------------------------------------------------------
$ws.Range("A13:C21").Clear()

$ws.Range("A100:C100").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A20:C20").PasteSpecial(-4122)

$ws.Range("A101:C101").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A21:C21").PasteSpecial(-4122)

$ws.Range("A102:C102").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Clean up the scratch template rows.
$ws.Range("A100:C102").Clear()

# --- values -----------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Programa resumido:"
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"

$ws.Cells.Item(14, 1).Value = "Short syllabus:"

$ws.Cells.Item(15, 1).Value = "Programa:"
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "01/01/2012"
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "01/01/2012"

$ws.Cells.Item(16, 1).Value = "Syllabus:"

$ws.Cells.Item(17, 1).Value = "Avaliação:"

$ws.Cells.Item(18, 1).Value = "Método:"
$ws.Cells.Item(18, 2).Value = "7455355 - Robson da Silva Rocha"
$ws.Cells.Item(18, 3).Value = "7455355 - Robson da Silva Rocha"

$ws.Cells.Item(19, 1).Value = "Critério:"
$ws.Cells.Item(19, 2).Value = "Serão ministradas aulas expositivas convencionais, associadas à exposição de vídeos e slides sobre sistemas de tratamento avançado. Além disso serão efetuadas visitas e serão desenvolvidos exercícios orientados."
$ws.Cells.Item(19, 3).Value = "Serão ministradas aulas expositivas convencionais, associadas à exposição de vídeos e slides sobre sistemas de tratamento avançado. Além disso serão efetuadas visitas e serão desenvolvidos exercícios orientados."

$ws.Cells.Item(20, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(20, 2).Value = "Serão aplicadas duas provas (1o. e 2o. bimestres), com peso 8,0 e também será entregue lista de exercícios, com peso 2,0."
$ws.Cells.Item(20, 3).Value = "Serão aplicadas duas provas (1o. e 2o. bimestres), com peso 8,0 e também será entregue lista de exercícios, com peso 2,0."

$ws.Cells.Item(21, 1).Value = "Bibliografia:"
$ws.Cells.Item(21, 2).Value = "Elaboração de monografia, com tema escolhido pelo docente, enfocando matéria em que o aluno demonstrou menor habilidade (peso: 3,0); e prova escrita sobre todfa a matéria da disciplina (peso: 7,0)."
$ws.Cells.Item(21, 3).Value = "Elaboração de monografia, com tema escolhido pelo docente, enfocando matéria em que o aluno demonstrou menor habilidade (peso: 3,0); e prova escrita sobre todfa a matéria da disciplina (peso: 7,0)."

# --- row heights --------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

Write-Host "Done."
